# Update the crossing-stress lookup table (Sheet1!A2:E36).
# The whole data block is rewritten in place since nearly every row's
# values shifted/changed, and four new rows were appended at the bottom.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("", 3, 30, "f", 1),
    @("", 3, 35, "f", 3),
    @("", 3, 99, "f", 3),
    @("", 3, 99, "t", 2),
    @("", 5, 25, "f", 2),
    @("", 5, 30, "f", 3),
    @("", 5, 35, "f", 3),
    @("", 5, 99, "f", 4),
    @("", 5, 25, "t", 1),
    @("", 5, 30, "t", 2),
    @("", 5, 35, "t", 3),
    @("", 5, 99, "t", 4),
    @("", 99, 25, "t", 3),
    @("", 99, 30, "t", 3),
    @("", 99, 99, "t", 4),
    @("", 99, 99, "f", 4),
    @("rrfb", 3, 30, "f", 1),
    @("rrfb", 3, 35, "f", 2),
    @("rrfb", 3, 99, "f", 3),
    @("rrfb", 3, 99, "t", 2),
    @("rrfb", 5, 25, "f", 2),
    @("rrfb", 5, 30, "f", 2),
    @("rrfb", 5, 35, "f", 3),
    @("rrfb", 5, 99, "f", 3),
    @("rrfb", 5, 25, "t", 1),
    @("rrfb", 5, 30, "t", 2),
    @("rrfb", 5, 35, "t", 2),
    @("rrfb", 5, 99, "t", 3),
    @("rrfb", 99, 25, "t", 3),
    @("rrfb", 99, 30, "t", 3),
    @("rrfb", 99, 99, "t", 4),
    @("rrfb", 99, 99, "f", 4),
    @("signal", 99, 99, "f", 1),
    @("hawk", 99, 99, "f", 1),
    @("four way stop", 99, 99, "f", 1)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r++
}

[void]$ws.Range("E22").Select()
